$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Preserve the cell as plain text even when $value looks like a number,
    # without leaving any lasting NumberFormat/style change behind.
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Cells.Item(2, 4).Value = '28.770.36'
$ws.Cells.Item(2, 5).Value = '  +2.30%  '

$ws.Cells.Item(3, 4).Value = '1.878.75'
$ws.Cells.Item(3, 5).Value = '  +2.51%  '

Set-TextValue $ws.Cells.Item(4, 4) '1.005'
$ws.Cells.Item(4, 5).Value = '  +0.33%  '

Set-TextValue $ws.Cells.Item(5, 4) '323.68'
$ws.Cells.Item(5, 5).Value = '  -1.71%  '

Set-TextValue $ws.Cells.Item(6, 4) '1.005'
$ws.Cells.Item(6, 5).Value = '  +0.35%  '

Set-TextValue $ws.Cells.Item(7, 4) '0.4673'
$ws.Cells.Item(7, 5).Value = '  +0.36%  '

Set-TextValue $ws.Cells.Item(8, 4) '0.3933'
$ws.Cells.Item(8, 5).Value = '  +1.49%  '

Set-TextValue $ws.Cells.Item(9, 4) '0.07920'
$ws.Cells.Item(9, 5).Value = '  +0.37%  '

$ws.Cells.Item(10, 5).Value = '  +2.12%  '

Set-TextValue $ws.Cells.Item(11, 4) '22.36'
$ws.Cells.Item(11, 5).Value = '  +1.65%  '

$ws.Cells.Item(12, 4).Value = '1.937.05'
$ws.Cells.Item(12, 5).Value = '  +4.77%  '

Set-TextValue $ws.Cells.Item(13, 4) '5.741'
$ws.Cells.Item(13, 5).Value = '  +1.20%  '

Set-TextValue $ws.Cells.Item(14, 4) '7.017'
$ws.Cells.Item(14, 5).Value = '  +1.58%  '

Set-TextValue $ws.Cells.Item(15, 4) '0.06978'
$ws.Cells.Item(15, 5).Value = '  +1.80%  '

Set-TextValue $ws.Cells.Item(16, 4) '88.70'
$ws.Cells.Item(16, 5).Value = '  +2.18%  '

Set-TextValue $ws.Cells.Item(17, 4) '1.006'
$ws.Cells.Item(17, 5).Value = '  +0.48%  '

Set-TextValue $ws.Cells.Item(18, 4) '0.00001009'
$ws.Cells.Item(18, 5).Value = '  +0.92%  '

Set-TextValue $ws.Cells.Item(19, 4) '16.96'
$ws.Cells.Item(19, 5).Value = '  +1.68%  '

$ws.Cells.Item(20, 5).Value = '  +0.35%  '

$ws.Cells.Item(21, 4).Value = '28.802.58'
$ws.Cells.Item(21, 5).Value = '  +2.33%  '

Set-TextValue $ws.Cells.Item(22, 4) '5.351'
$ws.Cells.Item(22, 5).Value = '  +0.27%  '

$ws.Cells.Item(23, 5).Value = '  +0.63%  '

Set-TextValue $ws.Cells.Item(24, 4) '2.127'
$ws.Cells.Item(24, 5).Value = '  +1.54%  '

$ws.Cells.Item(25, 4).Value = '2.159.17'
$ws.Cells.Item(25, 5).Value = '  +4.92%  '

Set-TextValue $ws.Cells.Item(26, 4) '153.32'
$ws.Cells.Item(26, 5).Value = '  +0.35%  '

Set-TextValue $ws.Cells.Item(27, 4) '19.36'
$ws.Cells.Item(27, 5).Value = '  +0.45%  '

Set-TextValue $ws.Cells.Item(28, 4) '5.767'
$ws.Cells.Item(28, 5).Value = '  -0.21%  '

Set-TextValue $ws.Cells.Item(29, 4) '2.000'
$ws.Cells.Item(29, 5).Value = '  +1.23%  '

Set-TextValue $ws.Cells.Item(30, 4) '119.73'
$ws.Cells.Item(30, 5).Value = '  +2.07%  '

Set-TextValue $ws.Cells.Item(31, 4) '0.09396'
$ws.Cells.Item(31, 5).Value = '  +1.32%  '

Set-TextValue $ws.Cells.Item(32, 4) '0.9378'
$ws.Cells.Item(32, 5).Value = '  -0.02%  '

Set-TextValue $ws.Cells.Item(33, 4) '5.311'
$ws.Cells.Item(33, 5).Value = '  +0.07%  '

Set-TextValue $ws.Cells.Item(34, 4) '1.355'
$ws.Cells.Item(34, 5).Value = '  +2.52%  '

$ws.Cells.Item(35, 5).Value = '  -0.32%  '

Set-TextValue $ws.Cells.Item(36, 4) '0.05908'
$ws.Cells.Item(36, 5).Value = '  -0.56%  '

Set-TextValue $ws.Cells.Item(37, 4) '0.02129'
$ws.Cells.Item(37, 5).Value = '  -1.10%  '

Set-TextValue $ws.Cells.Item(38, 4) '1.163'
$ws.Cells.Item(38, 5).Value = '  +1.07%  '

Set-TextValue $ws.Cells.Item(39, 4) '7.886'
$ws.Cells.Item(39, 5).Value = '  +2.54%  '

Set-TextValue $ws.Cells.Item(40, 4) '0.5724'
$ws.Cells.Item(40, 5).Value = '  +2.09%  '

$ws.Cells.Item(41, 5).Value = '  +1.22%  '

Set-TextValue $ws.Cells.Item(42, 4) '9.996'
$ws.Cells.Item(42, 5).Value = '  +0.59%  '

Set-TextValue $ws.Cells.Item(43, 4) '0.07311'
$ws.Cells.Item(43, 5).Value = '  +3.68%  '

Set-TextValue $ws.Cells.Item(44, 4) '11.83'
$ws.Cells.Item(44, 5).Value = '  +1.80%  '

Set-TextValue $ws.Cells.Item(45, 4) '1.180'
$ws.Cells.Item(45, 5).Value = '  -4.53%  '

Set-TextValue $ws.Cells.Item(46, 4) '0.5355'
$ws.Cells.Item(46, 5).Value = '  +1.37%  '

Set-TextValue $ws.Cells.Item(47, 4) '1.846'
$ws.Cells.Item(47, 5).Value = '  +0.32%  '

Set-TextValue $ws.Cells.Item(48, 4) '113.84'
$ws.Cells.Item(48, 5).Value = '  +1.57%  '

Set-TextValue $ws.Cells.Item(49, 4) '2.072'

Set-TextValue $ws.Cells.Item(50, 4) '2.373'
$ws.Cells.Item(50, 5).Value = '  +2.26%  '

Set-TextValue $ws.Cells.Item(51, 4) '1.004'
$ws.Cells.Item(51, 5).Value = '  +0.38%  '
